# Update the cached "last updated" date field text wherever it was
# hard-coded (slide master + every slide layout), and rename the
# "Users" box on the diagram slide to "Developers".

$p = $ppt.ActivePresentation

$oldDate = "24.07.2018"
$newDate = "01.08.2018"

function Update-DatePlaceholder($shapes) {
    for ($k = 1; $k -le $shapes.Count; $k++) {
        $shp = $shapes.Item($k)
        if ($shp.Name -like "Date Placeholder*") {
            if ($shp.TextFrame.TextRange.Text -eq $oldDate) {
                $shp.TextFrame.TextRange.Text = $newDate
            }
        }
    }
}

# Slide master.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout off the master.
$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholder $layouts.Item($i).Shapes
}

# Rename the "Users" rectangle to "Developers" on slide 1.
$s = $p.Slides.Item(1)
for ($j = 1; $j -le $s.Shapes.Count; $j++) {
    $shp = $s.Shapes.Item($j)
    if ($shp.HasTextFrame -and $shp.TextFrame.TextRange.Text -eq "Users") {
        $shp.TextFrame.TextRange.Text = "Developers"
    }
}
